$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4879.727
$ws.Range("I76").Value = 4964.3335
$ws.Range("J76").Value = 4499
$ws.Range("K76").Value = 4964.3335
$ws.Range("L76").Value = 4499
$ws.Range("M76").Value = -4649.3335
$ws.Range("N76").Value = -5129

$ws.Range("H79").Value = 4879.727
$ws.Range("I79").Value = 4964.3335
$ws.Range("J79").Value = 4499
$ws.Range("K79").Value = 4964.3335
$ws.Range("L79").Value = 4499
$ws.Range("M79").Value = -3872.3335
$ws.Range("N79").Value = -6683

$ws.Range("H92").Value = 2512.963
$ws.Range("I92").Value = 2607.6924
$ws.Range("K92").Value = 2607.6924
$ws.Range("M92").Value = -1359.6924

$ws.Range("H138").Value = 11194.118
$ws.Range("J138").Value = 11362.726
$ws.Range("L138").Value = 34088.178
$ws.Range("N138").Value = -44368.178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8425.478999999999
$ws.Range("I32").Value = 6763
$ws.Range("K32").Value = 6763
$ws.Range("M32").Value = -6476

$ws.Range("H61").Value = 6788.357
$ws.Range("J61").Value = 9466
$ws.Range("L61").Value = 9466
$ws.Range("N61").Value = -9890

$ws.Range("H74").Value = 4456.8945
$ws.Range("I74").Value = 4324.9375
$ws.Range("J74").Value = 5160.6665
$ws.Range("K74").Value = 4324.9375
$ws.Range("L74").Value = 5160.6665
$ws.Range("M74").Value = -3450.9375
$ws.Range("N74").Value = -6908.6665

$ws.Range("H77").Value = 4456.8945
$ws.Range("I77").Value = 4324.9375
$ws.Range("J77").Value = 5160.6665
$ws.Range("K77").Value = 21624.6875
$ws.Range("L77").Value = 25803.3325
$ws.Range("M77").Value = -17256.6875
$ws.Range("N77").Value = -34539.3325

$ws.Range("H132").Value = 3196.139
$ws.Range("I132").Value = 3108.0312
$ws.Range("J132").Value = 3901
$ws.Range("K132").Value = 9324.0936
$ws.Range("L132").Value = 11703
$ws.Range("M132").Value = -6794.0936
$ws.Range("N132").Value = -16763

$ws.Range("H136").Value = 6788.357
$ws.Range("J136").Value = 9466
$ws.Range("L136").Value = 28398
$ws.Range("N136").Value = -33498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 14168103
$ws.Range("I86").Value = 20239506
$ws.Range("J86").Value = 1494.3334
$ws.Range("K86").Value = 20239506
$ws.Range("L86").Value = 1494.3334
$ws.Range("M86").Value = -20238383
$ws.Range("N86").Value = -3740.3334

$ws.Range("H89").Value = 14168103
$ws.Range("I89").Value = 20239506
$ws.Range("J89").Value = 1494.3334
$ws.Range("K89").Value = 101197530
$ws.Range("L89").Value = 7471.666999999999
$ws.Range("M89").Value = -101191914
$ws.Range("N89").Value = -18703.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6225.737
$ws.Range("I31").Value = 5893.154
$ws.Range("J31").Value = 6946.3335
$ws.Range("K31").Value = 5893.154
$ws.Range("L31").Value = 6946.3335
$ws.Range("M31").Value = -5598.154
$ws.Range("N31").Value = -7536.3335

$ws.Range("H34").Value = 6225.737
$ws.Range("I34").Value = 5893.154
$ws.Range("J34").Value = 6946.3335
$ws.Range("K34").Value = 5893.154
$ws.Range("L34").Value = 6946.3335
$ws.Range("M34").Value = -5691.154
$ws.Range("N34").Value = -7350.3335

$ws.Range("H58").Value = 3306.7
$ws.Range("I58").Value = 3033.2812
$ws.Range("K58").Value = 3033.2812
$ws.Range("M58").Value = -2830.2812

$ws.Range("H99").Value = 4711.067
$ws.Range("I99").Value = 4691.2
$ws.Range("K99").Value = 4691.2
$ws.Range("M99").Value = -3193.2

$ws.Range("H126").Value = 4711.067
$ws.Range("I126").Value = 4691.2
$ws.Range("K126").Value = 14073.6
$ws.Range("M126").Value = -11603.6

$ws.Range("H134").Value = 1382.1765
$ws.Range("I134").Value = 1406.0625
$ws.Range("K134").Value = 4218.1875
$ws.Range("M134").Value = -1683.1875

$ws.Range("H136").Value = 3306.7
$ws.Range("I136").Value = 3033.2812
$ws.Range("K136").Value = 9099.8436
$ws.Range("M136").Value = -6549.8436

$ws.Range("H138").Value = 86801.164
$ws.Range("J138").Value = 86801.164
$ws.Range("L138").Value = 86801.164
$ws.Range("N138").Value = -97081.164

$ws.Range("H141").Value = 290791.56
$ws.Range("J141").Value = 290791.56
$ws.Range("L141").Value = 290791.56
$ws.Range("N141").Value = -301151.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1969.1818
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1969.1818
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 5907.5454
$ws.Range("N5").Value = -6131.5454
$ws.Range("M5").ClearContents()

$ws.Range("H112").Value = 9605.333000000001
$ws.Range("I112").Value = 7363
$ws.Range("J112").Value = 11399.2
$ws.Range("K112").Value = 22089
$ws.Range("L112").Value = 34197.60000000001
$ws.Range("M112").Value = -20981
$ws.Range("N112").Value = -36413.60000000001

$ws.Range("H122").Value = 2434.3333
$ws.Range("I122").Value = 2304
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 20736
$ws.Range("L122").Value = 22495.5
$ws.Range("N122").Value = -27395.5
$ws.Range("M122").Value = -18286

$ws.Range("H131").Value = 4167.7144
$ws.Range("I131").Value = 2498.5
$ws.Range("J131").Value = 4835.4
$ws.Range("K131").Value = 7495.5
$ws.Range("L131").Value = 14506.2
$ws.Range("M131").Value = -2455.5
$ws.Range("N131").Value = -24586.2

$ws.Range("H132").Value = 2671.1428
$ws.Range("J132").Value = 3080
$ws.Range("L132").Value = 27720
$ws.Range("N132").Value = -32780

$ws.Range("H135").Value = 1969.1818
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 1969.1818
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 17722.6362
$ws.Range("N135").Value = -22792.6362
$ws.Range("M135").ClearContents()

$ws.Range("H139").Value = 4773.2173
$ws.Range("I139").Value = 3852.5
$ws.Range("J139").Value = 5264.2666
$ws.Range("K139").Value = 11557.5
$ws.Range("L139").Value = 15792.7998
$ws.Range("M139").Value = -6417.5
$ws.Range("N139").Value = -26072.7998

$ws.Range("H140").Value = 5778.2856
$ws.Range("I140").Value = 2650
$ws.Range("K140").Value = 7950
$ws.Range("M140").Value = -2770

$ws.Range("H141").Value = 18380
$ws.Range("I141").Value = 12966.667
$ws.Range("K141").Value = 38900.001
$ws.Range("M141").Value = -33720.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 20809962
$ws.Range("I11").Value = 29877376
$ws.Range("J11").Value = 16780000
$ws.Range("K11").Value = 29877376
$ws.Range("L11").Value = 16780000
$ws.Range("M11").Value = -29877237
$ws.Range("N11").Value = -16780278

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15210.154
$ws.Range("I40").Value = 14702.444
$ws.Range("K40").Value = 14702.444
$ws.Range("M40").Value = -14566.444

$ws.Range("H93").Value = 2242.111
$ws.Range("I93").Value = 2022.5
$ws.Range("K93").Value = 2022.5
$ws.Range("M93").Value = -774.5

$ws.Range("H122").Value = 5466.8687
$ws.Range("I122").Value = 5272.8438
$ws.Range("J122").Value = 6501.6665
$ws.Range("K122").Value = 15818.5314
$ws.Range("L122").Value = 19504.9995
$ws.Range("M122").Value = -13368.5314
$ws.Range("N122").Value = -24404.9995

$ws.Range("H136").Value = 9027.541999999999
$ws.Range("I136").Value = 7442.154
$ws.Range("J136").Value = 10901.182
$ws.Range("K136").Value = 22326.462
$ws.Range("L136").Value = 32703.546
$ws.Range("M136").Value = -19776.462
$ws.Range("N136").Value = -37803.546

$ws.Range("H139").Value = 89999.2
$ws.Range("J139").Value = 89999.2
$ws.Range("L139").Value = 89999.2
$ws.Range("N139").Value = -100279.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2415
$ws.Range("I122").Value = 2415
$ws.Range("K122").Value = 7245
$ws.Range("M122").Value = -4795

$ws.Range("H126").Value = 1852.76
$ws.Range("I126").Value = 1679.9584
$ws.Range("K126").Value = 5039.8752
$ws.Range("M126").Value = -2569.8752

$ws.Range("H132").Value = 5828.16
$ws.Range("J132").Value = 8736.25
$ws.Range("L132").Value = 26208.75
$ws.Range("N132").Value = -31268.75
